# Scheduled runner update: refresh market-price-derived columns
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
#  LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ) across the
# per-job "Spriggan_Profits" sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

function Set-Cells {
    param(
        $ws,
        $updates
    )
    foreach ($u in $updates) {
        $ws.Range($u[0]).Value = $u[1]
    }
}

# ---------------------------------------------------------------- ALC ----
$ws = $wb.Worksheets.Item("ALC")
Set-Cells $ws @(
    ,@("H33", 443.8125)
    ,@("I33", 463.4)
    ,@("K33", 463.4)
    ,@("M33", -234.4)

    ,@("H76", 8661.727999999999)
    ,@("I76", 8282.286)
    ,@("K76", 8282.286)
    ,@("M76", -7967.286)

    ,@("H79", 8661.727999999999)
    ,@("I79", 8282.286)
    ,@("K79", 8282.286)
    ,@("M79", -7190.286)

    ,@("H106", 2933.4707)
    ,@("I106", 2182.4614)
    ,@("K106", 2182.4614)
    ,@("M106", -1551.4614)

    ,@("H137", 4135.074)
    ,@("J137", 4717.909)
    ,@("L137", 14153.727)
    ,@("N137", -19253.727)
)

# ---------------------------------------------------------------- ARM ----
$ws = $wb.Worksheets.Item("ARM")
Set-Cells $ws @(
    ,@("H61", 76941890)
    ,@("I61", 76941890)
    ,@("K61", 76941890)
    ,@("M61", -76941678)

    ,@("H74", 62506348)
    ,@("I74", 71435330)
    ,@("J74", 3500)
    ,@("K74", 71435330)
    ,@("L74", 3500)
    ,@("M74", -71434456)
    ,@("N74", -5248)

    ,@("H77", 62506348)
    ,@("I77", 71435330)
    ,@("J77", 3500)
    ,@("K77", 357176650)
    ,@("L77", 17500)
    ,@("M77", -357172282)
    ,@("N77", -26236)

    ,@("H102", 90909090)
    ,@("I102", 90909090)
    ,@("K102", 90909090)
    ,@("M102", -90907468)

    ,@("H136", 76941890)
    ,@("I136", 76941890)
    ,@("K136", 230825670)
    ,@("M136", -230823120)
)

# ---------------------------------------------------------------- BSM ----
$ws = $wb.Worksheets.Item("BSM")
Set-Cells $ws @(
    ,@("H19", 0)
    ,@("J19", 0)
    ,@("L19", 0)
)
# Row 19 no longer carries a LeveProfitHQ figure once LevePriceHQ is zero.
$ws.Range("N19").ClearContents()

# ---------------------------------------------------------------- CRP ----
$ws = $wb.Worksheets.Item("CRP")
Set-Cells $ws @(
    ,@("H16", 2175210.8)
    ,@("I16", 2718763.5)
    ,@("K16", 2718763.5)
    ,@("M16", -2718476.5)

    ,@("H31", 13846.936)
    ,@("J31", 17521.875)
    ,@("L31", 17521.875)
    ,@("N31", -18111.875)

    ,@("H34", 13846.936)
    ,@("J34", 17521.875)
    ,@("L34", 17521.875)
    ,@("N34", -17925.875)

    ,@("H86", 4264.5)
    ,@("I86", 4393.8)
    ,@("J86", 4135.2)
    ,@("K86", 4393.8)
    ,@("L86", 4135.2)
    ,@("M86", -3270.8)
    ,@("N86", -6381.2)

    ,@("H89", 4264.5)
    ,@("I89", 4393.8)
    ,@("J89", 4135.2)
    ,@("K89", 21969)
    ,@("L89", 20676)
    ,@("M89", -16353)
    ,@("N89", -31908)

    ,@("H99", 1825)
    ,@("I99", 1825)
    ,@("K99", 1825)
    ,@("M99", -327)

    ,@("H105", 4001019.5)
    ,@("I105", 5000899.5)
    ,@("K105", 5000899.5)
    ,@("M105", -4999152.5)

    ,@("H109", 57927)
    ,@("J109", 71890.5)
    ,@("L109", 71890.5)
    ,@("N109", -73970.5)

    ,@("H113", 2175210.8)
    ,@("I113", 2718763.5)
    ,@("K113", 2718763.5)
    ,@("M113", -2716593.5)

    ,@("H126", 1825)
    ,@("I126", 1825)
    ,@("K126", 5475)
    ,@("M126", -3005)
)

# ---------------------------------------------------------------- CUL ----
$ws = $wb.Worksheets.Item("CUL")
Set-Cells $ws @(
    ,@("H61", 560.25)
    ,@("I61", 150)
    ,@("J61", 697)
    ,@("K61", 450)
    ,@("L61", 2091)
    ,@("M61", -235)
    ,@("N61", -2521)

    ,@("H80", 3978.8)
    ,@("I80", 3989)
    ,@("J80", 3977.6667)
    ,@("K80", 11967)
    ,@("L80", 11933.0001)
    ,@("M80", -11031)
    ,@("N80", -13805.0001)

    ,@("H83", 3978.8)
    ,@("I83", 3989)
    ,@("J83", 3977.6667)
    ,@("K83", 35901)
    ,@("L83", 35799.0003)
    ,@("M83", -31221)
    ,@("N83", -45159.0003)

    ,@("H131", 1223.6522)
    ,@("I131", 813.8889)
    ,@("J131", 2698.8)
    ,@("K131", 2441.6667)
    ,@("L131", 8096.400000000001)
    ,@("M131", 2598.3333)
    ,@("N131", -18176.4)

    ,@("H132", 1512.0625)
    ,@("I132", 1493.3334)
    ,@("K132", 13440.0006)
    ,@("M132", -10910.0006)

    ,@("H140", 2640.8333)
    ,@("J140", 6824.5)
    ,@("L140", 20473.5)
    ,@("N140", -30833.5)
)

# ---------------------------------------------------------------- GSM ----
$ws = $wb.Worksheets.Item("GSM")
Set-Cells $ws @(
    ,@("H64", 59988.89)
    ,@("J64", 59988.89)
    ,@("L64", 59988.89)
    ,@("N64", -60484.89)

    ,@("H67", 59988.89)
    ,@("J67", 59988.89)
    ,@("L67", 59988.89)
    ,@("N67", -61704.89)

    ,@("H113", 74419.86)
    ,@("I113", 113334.664)
    ,@("J113", 4373.2)
    ,@("K113", 113334.664)
    ,@("L113", 4373.2)
    ,@("M113", -111164.664)
    ,@("N113", -8713.200000000001)

    ,@("H132", 4634996.5)
    ,@("I132", 5438952)
    ,@("K132", 16316856)
    ,@("M132", -16314326)

    ,@("H141", 85828.5)
    ,@("J141", 85828.5)
    ,@("L141", 85828.5)
    ,@("N141", -96188.5)
)

# ---------------------------------------------------------------- LTW ----
$ws = $wb.Worksheets.Item("LTW")
Set-Cells $ws @(
    ,@("H22", 5387.75)
    ,@("I22", 3850.3333)
    ,@("K22", 3850.3333)
    ,@("M22", -3555.3333)

    ,@("H27", 5387.75)
    ,@("I27", 3850.3333)
    ,@("K27", 3850.3333)
    ,@("M27", -3743.3333)

    ,@("H61", 2461.75)
    ,@("I61", 2327.375)
    ,@("K61", 2327.375)
    ,@("M61", -2125.375)

    ,@("H113", 2461.75)
    ,@("I113", 2327.375)
    ,@("K113", 2327.375)
    ,@("M113", -157.375)

    ,@("H133", 53584.75)
    ,@("J133", 53584.75)
    ,@("L133", 53584.75)
    ,@("N133", -58644.75)

    ,@("H136", 2904.4)
    ,@("I136", 2904.4)
    ,@("K136", 8713.200000000001)
    ,@("M136", -6163.200000000001)
)

# ---------------------------------------------------------------- WVR ----
$ws = $wb.Worksheets.Item("WVR")
Set-Cells $ws @(
    ,@("H62", 15900)
    ,@("I62", 4500)
    ,@("K62", 4500)
    ,@("M62", -3876)

    ,@("H65", 15900)
    ,@("I65", 4500)
    ,@("K65", 22500)
    ,@("M65", -19380)

    ,@("H96", 2489.3)
    ,@("J96", 2848.5)
    ,@("L96", 2848.5)
    ,@("N96", -5594.5)

    ,@("H107", 1289.2142)
    ,@("I107", 982.375)
    ,@("J107", 1698.3334)
    ,@("K107", 2947.125)
    ,@("L107", 5095.0002)
    ,@("M107", -1027.125)
    ,@("N107", -8935.0002)
)
